# GoodInfo_v2 - 2022-01-12 未完成
#
# 1) Row 11 (2022-01-11): B11/C11/E11 were stored as text ("56308.0", "-1656.0",
#    "-3264.0"); convert them to real numeric values (56308, -1656, -3264).
# 2) Append a new row 12 for 2022-01-12 with the same text-shaped values the
#    sheet previously used for "not yet final" rows (i.e. numbers stored as
#    text, e.g. "56308.0", "-935.0", "-1.66%", "0").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: B11, C11, E11 -> numeric -------------------------------------
$ws.Range("B11").Value = 56308
$ws.Range("C11").Value = -1656
$ws.Range("E11").Value = -3264

# --- Step 2: new row 12, values kept as text -------------------------------
# Force text storage (so "56308.0" etc. aren't auto-converted to numbers/
# dates/percentages by Excel's type inference) by temporarily applying a
# text number format, then restore the default "Normal" style so no extra
# formatting is left behind on the cells.
$textCells = @("A12", "B12", "C12", "D12", "E12")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A12").Value = "2022-01-12"
$ws.Range("B12").Value = "56308.0"
$ws.Range("C12").Value = "-935.0"
$ws.Range("D12").Value = "-1.66%"
$ws.Range("E12").Value = "0"

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# F12 / G12 stay blank (mirrors F11/G11 which are empty cells too); touch a
# no-op formatting property so the otherwise-untouched cells are still
# materialized in the sheet, matching the original row's column span.
$ws.Range("F12").Font.Bold = $false
$ws.Range("G12").Font.Bold = $false
